$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export rolling window moved forward by two days:
#   - drop the two oldest rows (2025-09-16, 2025-09-17)
#   - everything else shifts up two rows
#   - append a brand-new trailing row for 2025-12-13 (no data yet -> zeros)
$ws.Rows("2:3").Delete()

$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "2025-12-13"
$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122)
$ws.Range("B88").Value = 0
$ws.Range("C88").Value = 0
